$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# --- Row 12 (Tags) ---------------------------------------------------
# Before: B12=nutrition  C12=' protocol '  D12='phenotyping '  E12=Plant  F12=' metadata '  G12='study '  H12=MIAPPE
# After : B12=nutrition  C12=phenotyping   D12=Plant            E12=' metadata '  F12=study   G12=MIAPPE   (H12 cleared)
$ws.Range("C12").Value = "phenotyping"
$ws.Range("D12").Value = "Plant"
$ws.Range("E12").Value = " metadata "
$ws.Range("F12").Value = "study"
$ws.Range("G12").Value = "MIAPPE"
$ws.Range("H12").ClearContents()

# --- Row 13 (Tags Term Accession Number) ------------------------------
# Before: E13 = 'http://purl.obolibrary.org/obo/NCIT_C14258'
# After : C13 = 'DPBO:1000224', D13 = 'NCIT:C14258', E13 cleared, F13 = 'NCIT:C63536'
$ws.Range("C13").Value = "DPBO:1000224"
$ws.Range("D13").Value = "NCIT:C14258"
$ws.Range("D13").WrapText = $true
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = "NCIT:C63536"

# --- Row 14 (Tags Term Source REF) ------------------------------------
# Before: E14 = 'NCIT'
# After : D14 cleared, E14 cleared
$ws.Range("E14").ClearContents()

# Row 13 shrinks back to single-line height now the long URL is gone
$ws.Rows.Item(13).RowHeight = 28.8

# Move the active selection like the author ended up doing
$ws.Range("G19").Select()
